$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_QRS_TTS")
$ws.Name = "CRF_QRS_TTS"
